# Add data for 2022-02-23 (carjacking-by-month-yoy-historical)
#
# 1. Rename the sheet tab from "Through 2022-02-14" to "Through 2022-02-15"
# 2. Update the "2022 (through 02-14)" column header (I1) to "2022 (through 02-15)"
# 3. Bump the February 2022 count (I3) from 65 to 69
# 4. Bump the Total 2022 count (I14) from 226 to 230

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "Through 2022-02-15"

# Update header text for the 2022 column
$ws.Range("I1").Value = "2022 (through 02-15)"

# Update February count
$ws.Range("I3").Value = 69

# Update Total count
$ws.Range("I14").Value = 230
